$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Insert a new column before the old "Total Cost" column (D) rename stays the
# same column; insert one new column after it (for "Amount") and rename the
# old "Total Cost" header to "Price".
$ws.Range("D1").Value = "Price"
$ws.Columns("E").Insert()
$ws.Range("E1").Value = "Amount"

# Give the new "Amount" header the same bold style as the other plain-text
# headers (Phone Number / Item / Description), not the currency style copied
# from the "Price" column during the column insert.
$ws.Range("E1").NumberFormat = "General"

# Fill the rest of row 1 out to column V with the same header styling
$ws.Range("G1:V1").Font.Bold = $true

# Update column widths: column E ("Amount") is narrow, column F
# ("Description") is wide.
$ws.Columns("E").ColumnWidth = 8.6640625
$ws.Columns("F").ColumnWidth = 38.6640625

# Selection / frozen pane marker moves to D1
$ws.Range("D1").Select()

# The "data" sheet becomes the active tab
$ws.Activate()
